{"js": "// Replace the date line and every \"NN\u00d7NN=NNNN\" answer in the multiplication\n// table with the updated values from the new day's worksheet.\nconst replacements = [\n  [\"2025-12-25 Thursday\", \"2025-12-26 Friday\"],\n  [\"59\u00d762=3658\", \"87\u00d713=1131\"],\n  [\"71\u00d759=4189\", \"72\u00d751=3672\"],\n  [\"40\u00d718=720\", \"20\u00d790=1800\"],\n  [\"80\u00d736=2880\", \"52\u00d713=676\"],\n  [\"58\u00d721=1218\", \"12\u00d728=336\"],\n  [\"65\u00d764=4160\", \"47\u00d744=2068\"],\n  [\"96\u00d764=6144\", \"53\u00d778=4134\"],\n  [\"19\u00d721=399\", \"22\u00d719=418\"],\n  [\"49\u00d779=3871\", \"55\u00d735=1925\"],\n  [\"70\u00d778=5460\", \"28\u00d791=2548\"],\n  [\"60\u00d780=4800\", \"27\u00d712=324\"],\n  [\"78\u00d714=1092\", \"95\u00d726=2470\"],\n  [\"49\u00d727=1323\", \"62\u00d716=992\"],\n  [\"40\u00d757=2280\", \"56\u00d762=3472\"],\n  [\"59\u00d789=5251\", \"31\u00d779=2449\"],\n  [\"74\u00d759=4366\", \"52\u00d729=1508\"],\n  [\"81\u00d726=2106\", \"41\u00d756=2296\"],\n  [\"85\u00d784=7140\", \"55\u00d729=1595\"],\n  [\"54\u00d765=3510\", \"75\u00d722=1650\"],\n  [\"84\u00d776=6384\", \"22\u00d725=550\"],\n  [\"41\u00d759=2419\", \"76\u00d711=836\"],\n  [\"80\u00d793=7440\", \"96\u00d757=5472\"],\n  [\"96\u00d744=4224\", \"20\u00d725=500\"],\n  [\"26\u00d754=1404\", \"60\u00d724=1440\"],\n  [\"67\u00d770=4690\", \"28\u00d798=2744\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"NN\u00d7NN=NNNN\" answer in the multiplication\n# table with the updated values from the new day's worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-25 Thursday\", \"2025-12-26 Friday\"),\n    @(\"59\u00d762=3658\", \"87\u00d713=1131\"),\n    @(\"71\u00d759=4189\", \"72\u00d751=3672\"),\n    @(\"40\u00d718=720\", \"20\u00d790=1800\"),\n    @(\"80\u00d736=2880\", \"52\u00d713=676\"),\n    @(\"58\u00d721=1218\", \"12\u00d728=336\"),\n    @(\"65\u00d764=4160\", \"47\u00d744=2068\"),\n    @(\"96\u00d764=6144\", \"53\u00d778=4134\"),\n    @(\"19\u00d721=399\", \"22\u00d719=418\"),\n    @(\"49\u00d779=3871\", \"55\u00d735=1925\"),\n    @(\"70\u00d778=5460\", \"28\u00d791=2548\"),\n    @(\"60\u00d780=4800\", \"27\u00d712=324\"),\n    @(\"78\u00d714=1092\", \"95\u00d726=2470\"),\n    @(\"49\u00d727=1323\", \"62\u00d716=992\"),\n    @(\"40\u00d757=2280\", \"56\u00d762=3472\"),\n    @(\"59\u00d789=5251\", \"31\u00d779=2449\"),\n    @(\"74\u00d759=4366\", \"52\u00d729=1508\"),\n    @(\"81\u00d726=2106\", \"41\u00d756=2296\"),\n    @(\"85\u00d784=7140\", \"55\u00d729=1595\"),\n    @(\"54\u00d765=3510\", \"75\u00d722=1650\"),\n    @(\"84\u00d776=6384\", \"22\u00d725=550\"),\n    @(\"41\u00d759=2419\", \"76\u00d711=836\"),\n    @(\"80\u00d793=7440\", \"96\u00d757=5472\"),\n    @(\"96\u00d744=4224\", \"20\u00d725=500\"),\n    @(\"26\u00d754=1404\", \"60\u00d724=1440\"),\n    @(\"67\u00d770=4690\", \"28\u00d798=2744\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
